# "added 4wk low sales check"
# Updates the forecast values (MyForecast, Inventory Coverage, Stockout Risk,
# Seasonality Index) on the "Forecast Comparison" sheet, and the derived
# totals on the "Summary" sheet, to reflect a refreshed forecast run that now
# also accounts for a 4-week low-sales check.

$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison -------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W10)
$ws1.Range("D2").Value = 42
$ws1.Range("H2").Value = 4.81
$ws1.Range("L2").Value = 1.08

# Row 3 (W11)
$ws1.Range("D3").Value = 42
$ws1.Range("H3").Value = 3.78
$ws1.Range("L3").Value = 1.18

# Row 4 (W12)
$ws1.Range("D4").Value = 43
$ws1.Range("H4").Value = 2.76
$ws1.Range("L4").Value = 0.87

# Row 5 (W13)
$ws1.Range("D5").Value = 43
$ws1.Range("H5").Value = 1.76
$ws1.Range("L5").Value = 1.04

# Row 6 (W14)
$ws1.Range("D6").Value = 43
$ws1.Range("H6").Value = 0.76
$ws1.Range("I6").Value = "Low"
$ws1.Range("L6").Value = 0.9

# Row 7 (W15)
$ws1.Range("D7").Value = 43
$ws1.Range("L7").Value = 1.19

# Row 8 (W16)
$ws1.Range("D8").Value = 43
$ws1.Range("L8").Value = 0.86

# Row 9 (W17)
$ws1.Range("D9").Value = 43
$ws1.Range("L9").Value = 1.06

# Row 10 (W18)
$ws1.Range("D10").Value = 43
$ws1.Range("L10").Value = 1.06

# Row 11 (W19)
$ws1.Range("D11").Value = 43
$ws1.Range("L11").Value = 1.1

# Row 12 (W20)
$ws1.Range("D12").Value = 44
$ws1.Range("L12").Value = 1.18

# Row 13 (W21)
$ws1.Range("D13").Value = 44
$ws1.Range("L13").Value = 1.01

# Row 14 (W22)
$ws1.Range("D14").Value = 44
$ws1.Range("L14").Value = 0.96

# Row 15 (W23)
$ws1.Range("D15").Value = 44
$ws1.Range("L15").Value = 0.96

# Row 16 (W24)
$ws1.Range("D16").Value = 44
$ws1.Range("L16").Value = 0.91

# Row 17 (W25)
$ws1.Range("D17").Value = 44
$ws1.Range("L17").Value = 0.84

# --- Sheet: Summary ---------------------------------------------------------
# Column B on this sheet stores everything as text (see "N/A" rows), so force
# these numeric-looking updates to stay text too, matching the existing
# inlineStr cell type instead of letting COM auto-coerce them to numbers.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'698"
$ws2.Range("B10").Value = "'345"
$ws2.Range("B11").Value = "'171"
$ws2.Range("B12").Value = "'45"
$ws2.Range("B14").Value = "'42"
